# Apply the cibmtr-reporting-ig update to the "Metadata" sheet.
# (xl/worksheets/sheet1.xml in the underlying OOXML; "Include from RxNorm"
# / sheet2.xml is unaffected in content - only shared-string indices shift,
# which is an OOXML-serialization detail, not a visible cell-value change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "0.1.7"

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date refreshed
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Contact details replaced: row 10 becomes the publisher contact, row 11
# becomes a second, named contact (previously both rows were placeholder
# "No display for ContactDetail" text).
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# A new "Jurisdiction" property row is inserted after Contact, pushing
# Description / Purpose / Copyright / Immutable down by one row each.
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "RxNorm codes for Treosulfan"

$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""

$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""

# Immutable row, now pushed down to row 16 (new row - give it the same
# bordered/wrapped style as the rest of the data rows by copying row 15's
# formatting, same as Excel does when you fill a table down one row).
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
